$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 858.3333
$ws.Cells.Item(2, 10).Value = 1237.5
$ws.Cells.Item(2, 12).Value = 1237.5
$ws.Cells.Item(2, 14).Value = -1463.5

$ws.Cells.Item(39, 8).Value = 137.375
$ws.Cells.Item(39, 9).Value = 137.375
$ws.Cells.Item(39, 11).Value = 412.125
$ws.Cells.Item(39, 13).Value = -116.125

$ws.Cells.Item(40, 8).Value = 3999.88
$ws.Cells.Item(40, 10).Value = 3999.88
$ws.Cells.Item(40, 12).Value = 3999.88
$ws.Cells.Item(40, 14).Value = -4349.88

$ws.Cells.Item(70, 8).Value = 2199.7144
$ws.Cells.Item(70, 10).Value = 2316
$ws.Cells.Item(70, 12).Value = 6948
$ws.Cells.Item(70, 14).Value = -7488

$ws.Cells.Item(73, 8).Value = 2199.7144
$ws.Cells.Item(73, 10).Value = 2316
$ws.Cells.Item(73, 12).Value = 6948
$ws.Cells.Item(73, 14).Value = -8820

$ws.Cells.Item(87, 8).Value = 137670.5
$ws.Cells.Item(87, 10).Value = 110354
$ws.Cells.Item(87, 12).Value = 110354
$ws.Cells.Item(87, 14).Value = -112850

$ws.Cells.Item(90, 8).Value = 137670.5
$ws.Cells.Item(90, 10).Value = 110354
$ws.Cells.Item(90, 12).Value = 331062
$ws.Cells.Item(90, 14).Value = -343542

$ws.Cells.Item(100, 8).Value = 3478.4736
$ws.Cells.Item(100, 9).Value = 2687.28
$ws.Cells.Item(100, 11).Value = 2687.28
$ws.Cells.Item(100, 13).Value = -2146.28

$ws.Cells.Item(107, 8).Value = 1801.375
$ws.Cells.Item(107, 9).Value = 1874.8
$ws.Cells.Item(107, 11).Value = 1874.8
$ws.Cells.Item(107, 13).Value = 45.20000000000005

$ws.Cells.Item(125, 8).Value = 6346.5
$ws.Cells.Item(125, 10).Value = 7465.8
$ws.Cells.Item(125, 12).Value = 67192.2
$ws.Cells.Item(125, 14).Value = -72112.2

$ws.Cells.Item(132, 8).Value = 1907.225
$ws.Cells.Item(132, 9).Value = 1666.4231
$ws.Cells.Item(132, 11).Value = 4999.2693
$ws.Cells.Item(132, 13).Value = -2469.2693

$ws.Cells.Item(133, 8).Value = 120000
$ws.Cells.Item(133, 10).Value = 120000
$ws.Cells.Item(133, 12).Value = 120000
$ws.Cells.Item(133, 14).Value = -130120

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4149.129
$ws.Cells.Item(45, 9).Value = 4617.846
$ws.Cells.Item(45, 10).Value = 3810.611
$ws.Cells.Item(45, 11).Value = 4617.846
$ws.Cells.Item(45, 12).Value = 3810.611
$ws.Cells.Item(45, 13).Value = -4240.846
$ws.Cells.Item(45, 14).Value = -4564.611

$ws.Cells.Item(88, 8).Value = 1102.5
$ws.Cells.Item(88, 10).Value = 1096.8334
$ws.Cells.Item(88, 12).Value = 1096.8334
$ws.Cells.Item(88, 14).Value = -1908.8334

$ws.Cells.Item(91, 8).Value = 1102.5
$ws.Cells.Item(91, 10).Value = 1096.8334
$ws.Cells.Item(91, 12).Value = 1096.8334
$ws.Cells.Item(91, 14).Value = -3904.8334

$ws.Cells.Item(109, 8).Value = 53459.75
$ws.Cells.Item(109, 10).Value = 53459.75
$ws.Cells.Item(109, 12).Value = 53459.75
$ws.Cells.Item(109, 14).Value = -56233.75

$ws.Cells.Item(122, 8).Value = 5160.552
$ws.Cells.Item(122, 9).Value = 3999.75
$ws.Cells.Item(122, 11).Value = 11999.25
$ws.Cells.Item(122, 13).Value = -9549.25

$ws.Cells.Item(132, 8).Value = 3097.3704
$ws.Cells.Item(132, 9).Value = 3071.0435
$ws.Cells.Item(132, 11).Value = 9213.130500000001
$ws.Cells.Item(132, 13).Value = -6683.130500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 1202
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 1202
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 1202
$ws.Cells.Item(64, 13).ClearContents()
$ws.Cells.Item(64, 14).Value = -1652

$ws.Cells.Item(67, 8).Value = 1202
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 1202
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 1202
$ws.Cells.Item(67, 13).ClearContents()
$ws.Cells.Item(67, 14).Value = -2762

$ws.Cells.Item(80, 8).Value = 541.4545000000001
$ws.Cells.Item(80, 9).Value = 367.33334
$ws.Cells.Item(80, 10).Value = 606.75
$ws.Cells.Item(80, 11).Value = 367.33334
$ws.Cells.Item(80, 12).Value = 606.75
$ws.Cells.Item(80, 13).Value = 630.66666
$ws.Cells.Item(80, 14).Value = -2602.75

$ws.Cells.Item(83, 8).Value = 541.4545000000001
$ws.Cells.Item(83, 9).Value = 367.33334
$ws.Cells.Item(83, 10).Value = 606.75
$ws.Cells.Item(83, 11).Value = 1836.6667
$ws.Cells.Item(83, 12).Value = 3033.75
$ws.Cells.Item(83, 13).Value = 3155.3333
$ws.Cells.Item(83, 14).Value = -13017.75

$ws.Cells.Item(134, 8).Value = 3509975.2
$ws.Cells.Item(134, 9).Value = 3922808
$ws.Cells.Item(134, 11).Value = 11768424
$ws.Cells.Item(134, 13).Value = -11765889

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4844.595
$ws.Cells.Item(31, 9).Value = 2033.6111
$ws.Cells.Item(31, 11).Value = 2033.6111
$ws.Cells.Item(31, 13).Value = -1738.6111

$ws.Cells.Item(34, 8).Value = 4844.595
$ws.Cells.Item(34, 9).Value = 2033.6111
$ws.Cells.Item(34, 11).Value = 2033.6111
$ws.Cells.Item(34, 13).Value = -1831.6111

$ws.Cells.Item(122, 8).Value = 3609.7144
$ws.Cells.Item(122, 9).Value = 2924.25
$ws.Cells.Item(122, 10).Value = 4523.6665
$ws.Cells.Item(122, 11).Value = 8772.75
$ws.Cells.Item(122, 12).Value = 13570.9995
$ws.Cells.Item(122, 13).Value = -6322.75
$ws.Cells.Item(122, 14).Value = -18470.9995

$ws.Cells.Item(134, 8).Value = 2775.3333
$ws.Cells.Item(134, 9).Value = 2200.75
$ws.Cells.Item(134, 10).Value = 3924.5
$ws.Cells.Item(134, 11).Value = 6602.25
$ws.Cells.Item(134, 12).Value = 11773.5
$ws.Cells.Item(134, 13).Value = -4067.25
$ws.Cells.Item(134, 14).Value = -16843.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 4228.5
$ws.Cells.Item(55, 9).Value = 1075
$ws.Cells.Item(55, 10).Value = 5016.875
$ws.Cells.Item(55, 11).Value = 3225
$ws.Cells.Item(55, 12).Value = 15050.625
$ws.Cells.Item(55, 13).Value = -3048
$ws.Cells.Item(55, 14).Value = -15404.625

$ws.Cells.Item(68, 8).Value = 917.8
$ws.Cells.Item(68, 9).Value = 897
$ws.Cells.Item(68, 10).Value = 923
$ws.Cells.Item(68, 11).Value = 2691
$ws.Cells.Item(68, 12).Value = 2769
$ws.Cells.Item(68, 13).Value = -1880
$ws.Cells.Item(68, 14).Value = -4391

$ws.Cells.Item(71, 8).Value = 917.8
$ws.Cells.Item(71, 9).Value = 897
$ws.Cells.Item(71, 10).Value = 923
$ws.Cells.Item(71, 11).Value = 8073
$ws.Cells.Item(71, 12).Value = 8307
$ws.Cells.Item(71, 13).Value = -4017
$ws.Cells.Item(71, 14).Value = -16419

$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 13).ClearContents()
$ws.Cells.Item(105, 14).ClearContents()

$ws.Cells.Item(120, 8).Value = 21461.309
$ws.Cells.Item(120, 9).Value = 13499.25
$ws.Cells.Item(120, 11).Value = 40497.75
$ws.Cells.Item(120, 13).Value = -35659.75

$ws.Cells.Item(121, 8).Value = 11618.8
$ws.Cells.Item(121, 10).Value = 18904.834
$ws.Cells.Item(121, 12).Value = 56714.50199999999
$ws.Cells.Item(121, 14).Value = -59334.50199999999

$ws.Cells.Item(139, 8).Value = 4000
$ws.Cells.Item(139, 9).Value = 4000
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 11).Value = 12000
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 13).Value = -6860
$ws.Cells.Item(139, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 757.13635
$ws.Cells.Item(97, 9).Value = 540.1053000000001
$ws.Cells.Item(97, 11).Value = 540.1053000000001
$ws.Cells.Item(97, 13).Value = -44.10530000000006

$ws.Cells.Item(122, 8).Value = 1492.4546
$ws.Cells.Item(122, 9).Value = 1202.4286
$ws.Cells.Item(122, 11).Value = 3607.2858
$ws.Cells.Item(122, 13).Value = -1157.2858

$ws.Cells.Item(132, 8).Value = 3157.4358
$ws.Cells.Item(132, 9).Value = 2680.276
$ws.Cells.Item(132, 10).Value = 4541.2
$ws.Cells.Item(132, 11).Value = 8040.828
$ws.Cells.Item(132, 12).Value = 13623.6
$ws.Cells.Item(132, 13).Value = -5510.828
$ws.Cells.Item(132, 14).Value = -18683.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1709.6666
$ws.Cells.Item(16, 9).Value = 1709.6666
$ws.Cells.Item(16, 11).Value = 1709.6666
$ws.Cells.Item(16, 13).Value = -1539.6666

$ws.Cells.Item(22, 8).Value = 5562.5
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 5562.5
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 5562.5
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).Value = -6152.5

$ws.Cells.Item(27, 8).Value = 5562.5
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 5562.5
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 5562.5
$ws.Cells.Item(27, 13).ClearContents()
$ws.Cells.Item(27, 14).Value = -5776.5

$ws.Cells.Item(40, 8).Value = 2858.7778
$ws.Cells.Item(40, 10).Value = 2831.6667
$ws.Cells.Item(40, 12).Value = 2831.6667
$ws.Cells.Item(40, 14).Value = -3103.6667

$ws.Cells.Item(46, 8).Value = 7064.3447
$ws.Cells.Item(46, 9).Value = 1900
$ws.Cells.Item(46, 10).Value = 7248.7856
$ws.Cells.Item(46, 11).Value = 1900
$ws.Cells.Item(46, 12).Value = 7248.7856
$ws.Cells.Item(46, 13).Value = -1712
$ws.Cells.Item(46, 14).Value = -7624.7856

$ws.Cells.Item(55, 8).Value = 5002
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 5002
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 5002
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -5348

$ws.Cells.Item(82, 8).Value = 2213.1155
$ws.Cells.Item(82, 9).Value = 1583.5883
$ws.Cells.Item(82, 11).Value = 1583.5883
$ws.Cells.Item(82, 13).Value = -1222.5883

$ws.Cells.Item(85, 8).Value = 2213.1155
$ws.Cells.Item(85, 9).Value = 1583.5883
$ws.Cells.Item(85, 11).Value = 1583.5883
$ws.Cells.Item(85, 13).Value = -335.5882999999999

$ws.Cells.Item(100, 8).Value = 1900
$ws.Cells.Item(100, 9).Value = 1800
$ws.Cells.Item(100, 11).Value = 1800
$ws.Cells.Item(100, 13).Value = -1259

$ws.Cells.Item(132, 8).Value = 3306.7144
$ws.Cells.Item(132, 9).Value = 2436
$ws.Cells.Item(132, 11).Value = 7308
$ws.Cells.Item(132, 13).Value = -4778

$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

$ws.Cells.Item(136, 8).Value = 6147.952
$ws.Cells.Item(136, 9).Value = 4138.273
$ws.Cells.Item(136, 11).Value = 12414.819
$ws.Cells.Item(136, 13).Value = -9864.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1995.7106
$ws.Cells.Item(132, 9).Value = 1783.0344
$ws.Cells.Item(132, 11).Value = 5349.1032
$ws.Cells.Item(132, 13).Value = -2819.1032
